# Update Name of Algo
# Applies updated KNN-imputed values to the "terrestrial_mammals" result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.855
$ws.Range("B21").Value = 9.019
$ws.Range("B23").Value = 8.154
$ws.Range("C24").Value = -12.015
$ws.Range("B25").Value = 5.922000000000001
$ws.Range("C28").Value = -12.874
$ws.Range("C36").Value = -12.732
$ws.Range("C45").Value = -13.519
$ws.Range("C48").Value = -11.205
$ws.Range("C49").Value = -13.082
$ws.Range("C52").Value = -11.131
$ws.Range("B53").Value = 5.281
$ws.Range("C53").Value = -10.908
$ws.Range("C54").Value = -13.409
$ws.Range("B57").Value = 5.028999999999999
$ws.Range("B59").Value = 4.708
$ws.Range("B69").Value = 5.726
$ws.Range("C70").Value = -11.052
$ws.Range("B79").Value = 5.577
$ws.Range("B83").Value = 5.165999999999999
$ws.Range("C86").Value = -13.893
$ws.Range("C87").Value = -13.345
$ws.Range("B93").Value = 5.131000000000001
$ws.Range("C101").Value = -12.26
